$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 121; existing rows 121.. shift down to 123..
$ws.Rows.Item(121).Insert()
$ws.Rows.Item(121).Insert()

# New row 121 data
$ws.Cells.Item(121,1).Value = 3
$ws.Cells.Item(121,2).Value = "Femacal de La Calera"
$ws.Cells.Item(121,3).Value = "Coquimbo"
$ws.Cells.Item(121,4).Value = 44511
$ws.Cells.Item(121,5).Value = 5
$ws.Cells.Item(121,6).Value = "Fruta"
$ws.Cells.Item(121,7).Value = 100101
$ws.Cells.Item(121,8).Value = "Berries"
$ws.Cells.Item(121,9).Value = 100112025
$ws.Cells.Item(121,10).Value = "Frutilla"
$ws.Cells.Item(121,11).Value = "Sin especificar"
$ws.Cells.Item(121,12).Value = "Especial"
$ws.Cells.Item(121,13).Value = 125
$ws.Cells.Item(121,14).Value = 6000
$ws.Cells.Item(121,15).Value = 6000
$ws.Cells.Item(121,16).Value = 6000
$ws.Cells.Item(121,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(121,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(121,19).Value = 857
$ws.Cells.Item(121,20).Value = 7

# New row 122 data
$ws.Cells.Item(122,1).Value = 3
$ws.Cells.Item(122,2).Value = "Femacal de La Calera"
$ws.Cells.Item(122,3).Value = "Coquimbo"
$ws.Cells.Item(122,4).Value = 44511
$ws.Cells.Item(122,5).Value = 5
$ws.Cells.Item(122,6).Value = "Fruta"
$ws.Cells.Item(122,7).Value = 100101
$ws.Cells.Item(122,8).Value = "Berries"
$ws.Cells.Item(122,9).Value = 100112025
$ws.Cells.Item(122,10).Value = "Frutilla"
$ws.Cells.Item(122,11).Value = "Sin especificar"
$ws.Cells.Item(122,12).Value = "Segunda"
$ws.Cells.Item(122,13).Value = 80
$ws.Cells.Item(122,14).Value = 4000
$ws.Cells.Item(122,15).Value = 4000
$ws.Cells.Item(122,16).Value = 4000
$ws.Cells.Item(122,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(122,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(122,19).Value = 571
$ws.Cells.Item(122,20).Value = 7
